$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new "custom accuracy" (2 decimal place) rounded figures
$ws.Range("B5").Value = 23.86
$ws.Range("C5").Value = 17.88
$ws.Range("D5").Value = 0.87
$ws.Range("E5").Value = 52.3
$ws.Range("F5").Value = 42.72
$ws.Range("G5").Value = 18.86
$ws.Range("H5").Value = 71.91
$ws.Range("I5").Value = 28.9
$ws.Range("J5").Value = 13.03
$ws.Range("K5").Value = 19.54
$ws.Range("L5").Value = 21.19
$ws.Range("M5").Value = 22.03
$ws.Range("N5").Value = 6.04
$ws.Range("O5").Value = 18.65
$ws.Range("P5").Value = 27.26
$ws.Range("Q5").Value = 15.32
$ws.Range("R5").Value = 0.22
$ws.Range("S5").Value = 0.57
$ws.Range("T5").Value = 278.63
$ws.Range("U5").Value = 52.54
$ws.Range("V5").Value = 17.52
$ws.Range("W5").Value = 35.43
$ws.Range("X5").Value = 18.84
$ws.Range("Y5").Value = 2.81
$ws.Range("Z5").Value = 35.56
$ws.Range("AA5").Value = 15.31
$ws.Range("AB5").Value = 13.36
$ws.Range("AC5").Value = 15.97
$ws.Range("AD5").Value = 22.04
$ws.Range("AE5").Value = 0.37
$ws.Range("AF5").Value = 65.53
$ws.Range("AG5").Value = 9.73
$ws.Range("AH5").Value = 21.86

# Narrow a handful of columns from width 8 to width 7 (raw XML units)
$ws.Columns.Item(10).ColumnWidth = 6.17   # J
$ws.Columns.Item(11).ColumnWidth = 6.17   # K
$ws.Columns.Item(17).ColumnWidth = 6.17   # Q
$ws.Columns.Item(27).ColumnWidth = 6.17   # AA
$ws.Columns.Item(28).ColumnWidth = 6.17   # AB
$ws.Columns.Item(29).ColumnWidth = 6.17   # AC

# Drop the now-redundant last data row (row 6), shrinking the sheet dimension to A1:AH5
$ws.Rows.Item(6).Delete()
